# Auto-generated edit script: updates Kraken_Profits price/profit columns (H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets per scheduled runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 12307.538
$ws.Range("I70").Value = 5000
$ws.Range("J70").Value = 12916.5
$ws.Range("K70").Value = 15000
$ws.Range("L70").Value = 38749.5
$ws.Range("M70").Value = -14730
$ws.Range("N70").Value = -39289.5
$ws.Range("H73").Value = 12307.538
$ws.Range("I73").Value = 5000
$ws.Range("J73").Value = 12916.5
$ws.Range("K73").Value = 15000
$ws.Range("L73").Value = 38749.5
$ws.Range("M73").Value = -14064
$ws.Range("N73").Value = -40621.5
$ws.Range("H80").Value = 8633.25
$ws.Range("I80").Value = 1146.6
$ws.Range("K80").Value = 3439.8
$ws.Range("M80").Value = -2441.8
$ws.Range("H83").Value = 8633.25
$ws.Range("I83").Value = 1146.6
$ws.Range("K83").Value = 10319.4
$ws.Range("M83").Value = -5327.4
$ws.Range("H125").Value = 2557.25
$ws.Range("I125").Value = 2449.5
$ws.Range("J125").Value = 2665
$ws.Range("K125").Value = 22045.5
$ws.Range("L125").Value = 23985
$ws.Range("M125").Value = -19585.5
$ws.Range("N125").Value = -28905
$ws.Range("H129").Value = 1673.5
$ws.Range("I129").Value = 564.6667
$ws.Range("J129").Value = 5000
$ws.Range("K129").Value = 1694.0001
$ws.Range("L129").Value = 15000
$ws.Range("M129").Value = 3305.9999
$ws.Range("N129").Value = -25000
$ws.Range("H131").Value = 230.66667
$ws.Range("I131").Value = 230.66667
$ws.Range("K131").Value = 692.00001
$ws.Range("M131").Value = 4347.99999
$ws.Range("H137").Value = 2001
$ws.Range("I137").Value = 2001
$ws.Range("K137").Value = 6003
$ws.Range("M137").Value = -3453

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 30033.334
$ws.Range("I28").Value = 30033.334
$ws.Range("K28").Value = 30033.334
$ws.Range("M28").Value = -29841.334
$ws.Range("H99").Value = 30033.334
$ws.Range("I99").Value = 30033.334
$ws.Range("K99").Value = 30033.334
$ws.Range("M99").Value = -27038.334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 1252.3334
$ws.Range("I25").Value = 1203.5
$ws.Range("J25").Value = 1350
$ws.Range("K25").Value = 1203.5
$ws.Range("L25").Value = 1350
$ws.Range("M25").Value = -968.5
$ws.Range("N25").Value = -1820
$ws.Range("H56").Value = 29473
$ws.Range("J56").Value = 29473
$ws.Range("L56").Value = 29473
$ws.Range("N56").Value = -30951

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 756.7273
$ws.Range("I16").Value = 756.7273
$ws.Range("K16").Value = 756.7273
$ws.Range("M16").Value = -469.7273
$ws.Range("H113").Value = 756.7273
$ws.Range("I113").Value = 756.7273
$ws.Range("K113").Value = 756.7273
$ws.Range("M113").Value = 1413.2727
$ws.Range("H122").Value = 1000
$ws.Range("I122").Value = 1000
$ws.Range("K122").Value = 3000
$ws.Range("M122").Value = -550

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 566.1667
$ws.Range("I5").Value = 479.4
$ws.Range("K5").Value = 1438.2
$ws.Range("M5").Value = -1326.2
$ws.Range("H17").Value = 56.733334
$ws.Range("I17").Value = 32.9
$ws.Range("J17").Value = 104.4
$ws.Range("K17").Value = 98.69999999999999
$ws.Range("L17").Value = 313.2
$ws.Range("M17").Value = 70.30000000000001
$ws.Range("N17").Value = -651.2
$ws.Range("H34").Value = 1379.9
$ws.Range("I34").Value = 324.25
$ws.Range("J34").Value = 2083.6667
$ws.Range("K34").Value = 972.75
$ws.Range("L34").Value = 6251.000100000001
$ws.Range("M34").Value = -888.75
$ws.Range("N34").Value = -6419.000100000001
$ws.Range("H50").Value = 4500
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").Value = $null
$ws.Range("H53").Value = 4500
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").Value = $null
$ws.Range("H55").Value = 2318.5454
$ws.Range("J55").Value = 3300.8
$ws.Range("L55").Value = 9902.400000000001
$ws.Range("N55").Value = -10256.4
$ws.Range("H97").Value = 531.8182
$ws.Range("I97").Value = 541
$ws.Range("J97").Value = 520.8
$ws.Range("K97").Value = 1623
$ws.Range("L97").Value = 1562.4
$ws.Range("M97").Value = -1127
$ws.Range("N97").Value = -2554.4
$ws.Range("H132").Value = 2035
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").Value = $null
$ws.Range("H135").Value = 566.1667
$ws.Range("I135").Value = 479.4
$ws.Range("K135").Value = 4314.599999999999
$ws.Range("M135").Value = -1779.599999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").Value = $null
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4540.1
$ws.Range("J7").Value = 5000
$ws.Range("L7").Value = 5000
$ws.Range("N7").Value = -5224
$ws.Range("H22").Value = 2829.3157
$ws.Range("I22").Value = 2521.5715
$ws.Range("J22").Value = 3008.8333
$ws.Range("K22").Value = 2521.5715
$ws.Range("L22").Value = 3008.8333
$ws.Range("M22").Value = -2226.5715
$ws.Range("N22").Value = -3598.8333
$ws.Range("H27").Value = 2829.3157
$ws.Range("I27").Value = 2521.5715
$ws.Range("J27").Value = 3008.8333
$ws.Range("K27").Value = 2521.5715
$ws.Range("L27").Value = 3008.8333
$ws.Range("M27").Value = -2414.5715
$ws.Range("N27").Value = -3222.8333
$ws.Range("H56").Value = 30020.5
$ws.Range("I56").Value = 30020.5
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 30020.5
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -29329.5
$ws.Range("N56").Value = $null
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").Value = $null
$ws.Range("H68").Value = 3600
$ws.Range("I68").Value = 3600
$ws.Range("K68").Value = 3600
$ws.Range("M68").Value = -2851
$ws.Range("H71").Value = 3600
$ws.Range("I71").Value = 3600
$ws.Range("K71").Value = 18000
$ws.Range("M71").Value = -14256
$ws.Range("H122").Value = 6000
$ws.Range("I122").Value = 6000
$ws.Range("K122").Value = 18000
$ws.Range("M122").Value = -15550
$ws.Range("H126").Value = 4540.1
$ws.Range("J126").Value = 5000
$ws.Range("L126").Value = 15000
$ws.Range("N126").Value = -19940

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 40000
$ws.Range("I58").Value = 40000
$ws.Range("K58").Value = 40000
$ws.Range("M58").Value = -39692
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").Value = $null
$ws.Range("H122").Value = 4985
$ws.Range("I122").Value = 4985
$ws.Range("K122").Value = 14955
$ws.Range("M122").Value = -12505
$ws.Range("H126").Value = 7047.125
$ws.Range("I126").Value = 6339.5713
$ws.Range("K126").Value = 19018.7139
$ws.Range("M126").Value = -16548.7139
